# =====================================================================
# Update ZBP_05_kontakt_s_lidmi.xlsx: add 4 new weekly survey periods
#   sheet "data"   -> new columns BB:BE (29.3.-4.4, 5.-11.4, 19.-25.4,
#                     26.4.-2.5. 2021) for rows 2-61
#   sheet "pocetR" -> new columns BA:BD (same 4 periods) for rows 2-22
# Also refresh the "aktualizace" (last-updated) date in the two title
# strings from 7. 4. 2021 to 11. 5. 2021, and clear a stray date-time
# number format that had leaked onto BA29/BA57 on the "data" sheet.
# =====================================================================

$wb = $excel.ActiveWorkbook
$wsData   = $wb.Worksheets.Item("data")
$wsPocetR = $wb.Worksheets.Item("pocetR")

$newHeaders = @("29. 3.–4. 4. 2021", "5.–11. 4. 2021", "19.–25. 4. 2021", "26. 4.–2. 5. 2021")

$data1 = @{
    2 = @(11.5,11.5,13.5,14)
    3 = @(5,6,8,8)
    4 = @(5.5,5.5,6.5,6.5)
    5 = @(14.5,15,17,17)
    6 = @(14.5,14.5,17.5,17.5)
    7 = @(7,7,8.5,9.5)
    8 = @(13.5,14.5,17,18)
    9 = @(10.5,10.5,12.5,12.5)
    10 = @(11.5,11,13,13)
    11 = @(15,15.5,17,17.5)
    12 = @(10.5,10.5,12,12.5)
    13 = @(9.5,9.5,13,13.5)
    14 = @(13,13,15,15.5)
    15 = @(10,10,12,12.5)
    16 = @(11,11.5,13.5,14)
    17 = @(11.5,11.5,12.5,13)
    18 = @(11.5,11.5,14.5,15)
    19 = @(12,12,15,15)
    20 = @(21,21,22,21.5)
    21 = @(11.5,11.5,14,17)
    22 = @(7.5,7,9.5,10.5)
    23 = @(6,6,7,8)
    24 = @(10,10,10,10)
    25 = @(8,8,10,10)
    26 = @(4,4,5,5)
    27 = @(8,9,10,11)
    28 = @(5,5,7,8)
    29 = @(5,5,6,6)
    30 = @(8,8,10,11)
    31 = @(5,5,7,7)
    32 = @(5,5,6,7)
    33 = @(7,7,10,10)
    34 = @(5,5,6,7)
    35 = @(5,5,8,9)
    36 = @(6,6,7,7)
    37 = @(6,6,8,10)
    38 = @(5,6,8,7)
    39 = @(15,15,15,15)
    40 = @(10,10,10,11)
    41 = @(5,5,6,7)
    42 = @(5,3,5,5)
    43 = @(7.5,7.5,10,10)
    44 = @(6,6,7,7.5)
    45 = @(3.5,3.5,4,4)
    46 = @(7,7,9.5,9)
    47 = @(5,5,5.5,5.5)
    48 = @(4.5,4.5,6,6)
    49 = @(6.5,6.5,8.5,8.5)
    50 = @(4.5,4.5,5.5,5.5)
    51 = @(5,5,6.5,6.5)
    52 = @(6,6,7,7)
    53 = @(4.5,4.5,6,6)
    54 = @(5.5,5.5,7,7)
    55 = @(5,5,5.5,5.5)
    56 = @(5.5,5.5,7,7)
    57 = @(4.5,4.5,6.5,6.5)
    58 = @(7.5,7.5,8.5,8.5)
    59 = @(5,5.5,7,8)
    60 = @(4.5,4.5,6,6)
    61 = @(4,4,4.5,5)
}

$data2 = @{
    2 = @(1832,1829,1815,1820)
    3 = @(421,419,414,413)
    4 = @(637,640,652,653)
    5 = @(774,770,749,754)
    6 = @(364,362,418,422)
    7 = @(640,639,710,711)
    8 = @(828,828,687,687)
    9 = @(496,491,445,448)
    10 = @(653,653,556,561)
    11 = @(683,685,814,811)
    12 = @(898,896,890,892)
    13 = @(934,933,925,928)
    14 = @(959,954,940,945)
    15 = @(410,411,411,414)
    16 = @(222,224,220,218)
    17 = @(241,240,244,243)
    18 = @(624,624,693,693)
    19 = @(85,85,88,87)
    20 = @(159,158,127,127)
    21 = @(121,123,76,77)
}
# ---- Sheet "data": header row (row 1), columns BB:BE = 54..57 ----
for ($i = 0; $i -lt 4; $i++) {
    $col = 54 + $i
    $srcCell = $wsData.Cells.Item(1, 53)   # BA1 -- existing header, same look
    $srcCell.Copy()
    $dstCell = $wsData.Cells.Item(1, $col)
    $dstCell.PasteSpecial(-4122)           # xlPasteFormats
    $dstCell.Value = $newHeaders[$i]
}
$wsData.Application.CutCopyMode = $false

# ---- Sheet "data": body rows 2-61, columns BB:BE = 54..57 ----
foreach ($r in $data1.Keys) {
    $vals = $data1[$r]
    for ($i = 0; $i -lt 4; $i++) {
        $wsData.Cells.Item($r, 54 + $i).Value2 = $vals[$i]
    }
}

# ---- Sheet "pocetR": header row (row 1), columns BA:BD = 53..56 ----
for ($i = 0; $i -lt 4; $i++) {
    $col = 53 + $i
    $srcCell = $wsPocetR.Cells.Item(1, 52)  # AZ1 -- existing header, same look
    $srcCell.Copy()
    $dstCell = $wsPocetR.Cells.Item(1, $col)
    $dstCell.PasteSpecial(-4122)            # xlPasteFormats
    $dstCell.Value = $newHeaders[$i]
}
$wsPocetR.Application.CutCopyMode = $false

# ---- Sheet "pocetR": body rows 2-21, columns BA:BD = 53..56 ----
foreach ($r in $data2.Keys) {
    $vals = $data2[$r]
    for ($i = 0; $i -lt 4; $i++) {
        $wsPocetR.Cells.Item($r, 53 + $i).Value2 = $vals[$i]
    }
}


# ---- Refresh "aktualizace" date in the two summary title strings ----
$wsData.Cells.Item(62, 1).Value = "Život během pandemie, Kontakt s lidmi, průměr celkově a ve skupinách, aktualizace 11. 5. 2021"
$wsPocetR.Cells.Item(22, 1).Value = "Život během pandemie, Kontakt s lidmi, velikost dotázaného souboru celkově a ve skupinách, aktualizace 11. 5. 2021"

# ---- Clear stray date-time number format that had leaked onto BA29/BA57 ----
$wsData.Range("BA29").Style = "Normal"
$wsData.Range("BA57").Style = "Normal"
